$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "/Users/David/Documents/projects/mastergui/"

# Strip the leading prefix from the PATH_HCP values in column H (rows 2-11;
# rows 7-11 duplicate the same five paths used in rows 2-6).
foreach ($r in 2..11) {
    $cell = $ws.Cells.Item($r, 8)
    $val = $cell.Value2
    if ($val -like "$prefix*") {
        $cell.Value2 = $val.Substring($prefix.Length)
    }
}

# Update the active selection from H4 to H3.
$ws.Range("H3").Select()
